$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel;
# force them to remain Text-typed (matching the source data, which is all text),
# then restore the default "Normal" style so no stray formatting is introduced.
$numericLookingCells = @(
    "D5", "D6", "D8", "D10", "D11", "D12", "D13", "D14", "D16", "D19", "D23", "D24", "D25", "D26",
    "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D38", "D39", "D41", "D43", "D45",
    "D46", "D47", "D49", "D51"
)
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cryptocurrency price/volume figures plus the two swapped rows (41 & 43)
$ws.Range("D2").Value = "69.053.18"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "3.472.70"
$ws.Range("E3").Value = "  -3.77%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "580.12"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "181.58"
$ws.Range("E6").Value = "  -4.95%  "
$ws.Range("D7").Value = "3.464.14"
$ws.Range("E7").Value = "  -3.86%  "
$ws.Range("D8").Value = "0.606"
$ws.Range("E8").Value = "  -4.16%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "0.196"
$ws.Range("E10").Value = "  +7.09%  "
$ws.Range("D11").Value = "0.638"
$ws.Range("E11").Value = "  -4.06%  "
$ws.Range("D12").Value = "53.26"
$ws.Range("E12").Value = "  -5.24%  "
$ws.Range("D13").Value = "0.0000300"
$ws.Range("E13").Value = "  -3.96%  "
$ws.Range("D14").Value = "9.33"
$ws.Range("E14").Value = "  -4.37%  "
$ws.Range("D15").Value = "4.032.31"
$ws.Range("E15").Value = "  -3.65%  "
$ws.Range("D16").Value = "19.11"
$ws.Range("E16").Value = "  -4.52%  "
$ws.Range("D17").Value = "69.105.57"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("D18").Value = "3.473.27"
$ws.Range("E18").Value = "  -3.71%  "
$ws.Range("D19").Value = "12.20"
$ws.Range("E19").Value = "  -4.22%  "
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("E21").Value = "  +9.86%  "
$ws.Range("E22").Value = "  -4.99%  "
$ws.Range("D23").Value = "18.52"
$ws.Range("E23").Value = "  -7.93%  "
$ws.Range("D24").Value = "4.49"
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("D25").Value = "4.83"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").Value = "94.58"
$ws.Range("E26").Value = "  -3.18%  "
$ws.Range("D27").Value = "11.02"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "2.94"
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("D29").Value = "9.01"
$ws.Range("E29").Value = "  -5.29%  "
$ws.Range("D30").Value = "31.42"
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("D31").Value = "7.14"
$ws.Range("E31").Value = "  -6.55%  "
$ws.Range("D32").Value = "12.39"
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("D33").Value = "63.55"
$ws.Range("E33").Value = "  -4.32%  "
$ws.Range("E34").Value = "  -6.25%  "
$ws.Range("D35").Value = "523.95"
$ws.Range("E35").Value = "  -9.61%  "
$ws.Range("D36").Value = "0.403"
$ws.Range("E36").Value = "  +0.71%  "
$ws.Range("E37").Value = "  +5.23%  "
$ws.Range("D38").Value = "37.64"
$ws.Range("E38").Value = "  -3.61%  "
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").Value = "0.0₃0750"
$ws.Range("E40").Value = "  -8.49%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "3.35"
$ws.Range("E41").Value = "  -3.23%  "
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "3.07"
$ws.Range("E43").Value = "  -6.86%  "
$ws.Range("D44").Value = "3.313.67"
$ws.Range("E44").Value = "  +2.74%  "
$ws.Range("D45").Value = "3.46"
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("D46").Value = "2.94"
$ws.Range("E46").Value = "  -4.49%  "
$ws.Range("D47").Value = "0.0436"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("E48").Value = "  -4.11%  "
$ws.Range("D49").Value = "8.85"
$ws.Range("E49").Value = "  -8.32%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "136.79"
$ws.Range("E51").Value = "  -1.04%  "

# Restore default styling on the cells we temporarily reformatted as Text
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
